# The deck ships with two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (Office blue palette)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" palette (the one
#                            actually applied to the slide master / slides)
#
# The authored change swaps the two themes' contents: theme2.xml (the
# theme that is actually applied to the deck) becomes the default
# "Office Theme" colour scheme, while theme1.xml becomes the "Integral"
# / "Red Violet" scheme. Net visible effect on the deck: every slide's
# applied theme switches from "Integral" (magenta/violet/blue accents)
# to the plain "Office Theme" (blue/orange/grey accents).
#
# PowerPoint's object model edits theme colours through
# ThemeColorScheme (DrawingML clrScheme slots dk1,lt1,dk2,lt2,
# accent1-6,hlink,folHlink -> ThemeColorSchemeIndex 1..12), reachable
# from a Slide/Design. Updating it updates the theme part backing the
# presentation's slide master (theme2.xml), matching the target diff.

function Get-RgbFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" colour scheme (dk1..folHlink), in
# ThemeColorSchemeIndex order (1=dk1, 2=lt1, 3=dk2, 4=lt2,
# 5=accent1 .. 10=accent6, 11=hlink, 12=folHlink).
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$themeColors = $s1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColor = $themeColors.Colors($i)
    $themeColor.RGB = Get-RgbFromHex $officeThemeColors[$i - 1]
}
